$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data in rows 15 and 16 have been swapped (all the "record" columns),
# the Ost/Nord (Q/R) coordinates are re-rounded to whole numbers, and the
# Starttid/Sluttid (Z/AB) time cells are cleared for both rows.

# --- Row 15 gets the values that used to belong to row 16 ---
$ws.Range("A15").Value = 112079439
$ws.Range("B15").Value = 90689
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 5966
$ws.Range("F15").Value = "Motaggsvamp"
$ws.Range("G15").Value = "Sarcodon squamosus"
$ws.Range("H15").Value = "(Schaeff.) Quél."
$ws.Range("Q15").Value = 563408
$ws.Range("R15").Value = 6576469
$ws.Range("Z15").ClearContents()
$ws.Range("AB15").ClearContents()

# --- Row 16 gets the values that used to belong to row 15 ---
$ws.Range("A16").Value = 112079417
$ws.Range("B16").Value = 8377
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 106545
$ws.Range("F16").Value = "Mindre märgborre"
$ws.Range("G16").Value = "Tomicus minor"
$ws.Range("H16").Value = "(Hartig, 1834)"
$ws.Range("Q16").Value = 563452
$ws.Range("R16").Value = 6576051
$ws.Range("Z16").ClearContents()
$ws.Range("AB16").ClearContents()
